# The deck's two theme parts (ppt/theme/theme1.xml and ppt/theme/theme2.xml)
# get swapped: theme1.xml (currently "Office Theme" / "Office" colors, used
# by the notes master) and theme2.xml (currently "Integral" / "Red Violet"
# colors, used by the slide master / presentation) exchange their color
# schemes.
#
# The PowerPoint object model only exposes one editable "theme" endpoint in
# this host: the (slide) Master's Theme.ThemeColorScheme, which is what gets
# written back out to ppt/theme/theme2.xml. We drive that object to hold the
# 12 "Office" theme colors that theme1.xml had before the edit, reproducing
# the color-scheme half of the swap for the reachable theme part.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Index -> scheme slot (PowerPoint's ThemeColorScheme ordering):
#  1 dk1   2 lt1   3 dk2   4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
# 11 hlink   12 folHlink
#
# Target values are the "Office" theme colors (RGB hex -> VBA BGR-packed
# long, i.e. R + G*256 + B*65536).
$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
